$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (the default/unstyled format used by data cells in columns B-E)
# captured up-front so it can be restored after forcing text interpretation
# on numeric-looking price strings (preventing Excel from auto-converting them
# to floating point numbers and losing exact decimal formatting).
$normalStyle = $ws.Range('B2').Style

$ws.Range('D2').Value = '64.918.71'
$ws.Range('E2').Value = '  +1.31%  '
$ws.Range('D3').Value = '3.392.23'
$ws.Range('E3').Value = '  +1.24%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = "'559.46"
$ws.Range('D5').Style = $normalStyle
$ws.Range('E5').Value = '  +1.53%  '
$ws.Range('D6').Value = "'175.02"
$ws.Range('D6').Style = $normalStyle
$ws.Range('E6').Value = '  +0.97%  '
$ws.Range('E7').Value = '  +1.80%  '
$ws.Range('D8').Value = '3.381.65'
$ws.Range('E8').Value = '  +1.23%  '
$ws.Range('E9').Value = '  +0.05%  '
$ws.Range('E10').Value = '  +11.31%  '
$ws.Range('D11').Value = "'0.630"
$ws.Range('D11').Style = $normalStyle
$ws.Range('E11').Value = '  +2.64%  '
$ws.Range('D12').Value = "'54.88"
$ws.Range('D12').Style = $normalStyle
$ws.Range('E12').Value = '  +1.73%  '
$ws.Range('D13').Value = "'0.0000280"
$ws.Range('D13').Style = $normalStyle
$ws.Range('E13').Value = '  +5.08%  '
$ws.Range('D14').Value = "'9.13"
$ws.Range('D14').Style = $normalStyle
$ws.Range('E14').Value = '  +2.52%  '
$ws.Range('D15').Value = '3.931.83'
$ws.Range('E15').Value = '  +5.44%  '
$ws.Range('D16').Value = "'18.35"
$ws.Range('D16').Style = $normalStyle
$ws.Range('E16').Value = '  +0.74%  '
$ws.Range('E17').Value = '  +1.67%  '
$ws.Range('D18').Value = '3.396.98'
$ws.Range('E18').Value = '  +1.68%  '
$ws.Range('D19').Value = '64.960.13'
$ws.Range('E19').Value = '  +1.63%  '
$ws.Range('D20').Value = "'11.84"
$ws.Range('D20').Style = $normalStyle
$ws.Range('E20').Value = '  +0.79%  '
$ws.Range('E21').Value = '  +1.42%  '
$ws.Range('D22').Value = "'475.34"
$ws.Range('D22').Style = $normalStyle
$ws.Range('E22').Value = '  +15.67%  '
$ws.Range('D23').Value = "'5.00"
$ws.Range('D23').Style = $normalStyle
$ws.Range('E23').Value = '  +13.64%  '
$ws.Range('E24').Value = '  +1.92%  '
$ws.Range('D25').Value = "'87.05"
$ws.Range('D25').Style = $normalStyle
$ws.Range('E25').Value = '  +4.82%  '
$ws.Range('D26').Value = "'13.49"
$ws.Range('D26').Style = $normalStyle
$ws.Range('E26').Value = '  -2.79%  '
$ws.Range('D27').Value = "'2.91"
$ws.Range('D27').Style = $normalStyle
$ws.Range('E27').Value = '  +6.36%  '
$ws.Range('D28').Value = "'10.90"
$ws.Range('D28').Style = $normalStyle
$ws.Range('E28').Value = '  +2.98%  '
$ws.Range('D29').Value = "'8.79"
$ws.Range('D29').Style = $normalStyle
$ws.Range('E29').Value = '  +1.56%  '
$ws.Range('D30').Value = "'31.12"
$ws.Range('D30').Style = $normalStyle
$ws.Range('E30').Value = '  +6.69%  '
$ws.Range('D31').Value = "'6.69"
$ws.Range('D31').Style = $normalStyle
$ws.Range('E31').Value = '  +4.40%  '
$ws.Range('D32').Value = "'11.55"
$ws.Range('D32').Style = $normalStyle
$ws.Range('E32').Value = '  +1.47%  '
$ws.Range('D33').Value = "'61.83"
$ws.Range('D33').Style = $normalStyle
$ws.Range('E33').Value = '  +6.36%  '
$ws.Range('D34').Value = "'571.51"
$ws.Range('D34').Style = $normalStyle
$ws.Range('E34').Value = '  -1.72%  '
$ws.Range('E35').Value = '  +1.33%  '
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('B37').Value = 'Stacks'
$ws.Range('C37').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D37').Value = "'3.55"
$ws.Range('D37').Style = $normalStyle
$ws.Range('E37').Value = '  +4.36%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').Value = "'0.140"
$ws.Range('D38').Style = $normalStyle
$ws.Range('E38').Value = '  -5.27%  '
$ws.Range('D39').Value = "'35.75"
$ws.Range('D39').Style = $normalStyle
$ws.Range('E39').Value = '  +1.35%  '
$ws.Range('D40').Value = '0.0₃0758'
$ws.Range('E40').Value = '  +2.29%  '
$ws.Range('D41').Value = "'0.372"
$ws.Range('D41').Style = $normalStyle
$ws.Range('E41').Value = '  +1.15%  '
$ws.Range('D42').Value = '3.092.54'
$ws.Range('E42').Value = '  -1.64%  '
$ws.Range('D43').Value = "'0.999"
$ws.Range('D43').Style = $normalStyle
$ws.Range('E43').Value = '  +0.21%  '
$ws.Range('D44').Value = "'2.86"
$ws.Range('D44').Style = $normalStyle
$ws.Range('E44').Value = '  +1.59%  '
$ws.Range('E45').Value = '  +3.40%  '
$ws.Range('E46').Value = '  +5.62%  '
$ws.Range('E47').Value = '  +1.83%  '
$ws.Range('D48').Value = "'3.14"
$ws.Range('D48').Style = $normalStyle
$ws.Range('E48').Value = '  -3.76%  '
$ws.Range('E49').Value = '  +0.01%  '
$ws.Range('D50').Value = "'137.70"
$ws.Range('D50').Style = $normalStyle
$ws.Range('E50').Value = '  +3.68%  '
$ws.Range('D51').Value = "'8.35"
$ws.Range('D51').Style = $normalStyle
$ws.Range('E51').Value = '  +3.13%  '
